$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.184.87"
$ws.Range("E2").Value = "  -4.14%  "

$ws.Range("D3").Value = "3.297.15"
$ws.Range("E3").Value = "  -4.45%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "559.16"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.42%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.45"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.08%  "

$ws.Range("E7").Value = "  +0.25%  "

$ws.Range("D8").Value = "3.294.78"
$ws.Range("E8").Value = "  -4.50%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.482"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.72%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.84"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.27%  "

$ws.Range("E11").Value = "  -4.07%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.406"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.40%  "

$ws.Range("D13").Value = "3.868.71"
$ws.Range("E13").Value = "  -4.29%  "

$ws.Range("E14").Value = "  +0.63%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "27.23"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.95%  "

$ws.Range("D16").Value = "3.282.44"
$ws.Range("E16").Value = "  -4.88%  "

$ws.Range("E17").Value = "  -4.18%  "

$ws.Range("D18").Value = "60.199.09"
$ws.Range("E18").Value = "  -4.17%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.11"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.58%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.31"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.24%  "

$ws.Range("E21").Value = "  -5.12%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "373.22"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.64%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "73.89"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.90%  "

$ws.Range("E24").Value = "  -3.57%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.02%  "

$ws.Range("D26").Value = "3.462.63"
$ws.Range("E26").Value = "  -3.42%  "

$ws.Range("E27").Value = "  -9.35%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.172"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.71%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.11%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.19"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -6.81%  "

$ws.Range("E31").Value = "  +0.01%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.04"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.50%  "

$ws.Range("E33").Value = "  -5.18%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "22.54"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.07%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.27"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.12%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.18"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.23%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "166.35"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.71%  "

$ws.Range("E38").Value = "  -2.92%  "

$ws.Range("E39").Value = "  -7.41%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "27.49"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -14.73%  "

$ws.Range("D41").Value = "3.332.86"
$ws.Range("E41").Value = "  -4.35%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0735"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.58%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "41.82"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.43%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.750"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.59%  "

$ws.Range("E45").Value = "  -4.14%  "

$ws.Range("E46").Value = "  -5.64%  "

$ws.Range("D48").Value = "2.365.89"
$ws.Range("E48").Value = "  -7.91%  "

$ws.Range("E49").Value = "  -0.04%  "

$ws.Range("E50").Value = "  -4.90%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "21.52"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.71%  "
